$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous table (old layout: A1:I7) before writing the new one.
$ws.Range("A1:O7").Clear()

# Header row (columns reordered + 6 new trailing columns appended).
$header = New-Object 'object[,]' 1,15
$header[0,0]  = "id"
$header[0,1]  = "name"
$header[0,2]  = "brand"
$header[0,3]  = "model"
$header[0,4]  = "dtyp_code"
$header[0,5]  = "min_driver_ver"
$header[0,6]  = "descr"
$header[0,7]  = "lang_code"
$header[0,8]  = "is_active"
$header[0,9]  = "cr_by"
$header[0,10] = "cr_dtimes"
$header[0,11] = "upd_by"
$header[0,12] = "upd_dtimes"
$header[0,13] = "is_deleted"
$header[0,14] = "del_dtimes"
$ws.Range("A1:O1").Value = $header

$crDtimes = 45079.578168229164

$data = New-Object 'object[,]' 6,15

# Row 2 (id 165 - Fingerprint Scanner / Aratek)
$data[0,0]  = 165
$data[0,1]  = "Fingerprint Scanner"
$data[0,2]  = "Aratek"
$data[0,3]  = "FRO900"
$data[0,4]  = "FRS"
$data[0,5]  = 1.12
$data[0,6]  = "Pour capturer les empreintes digitales"
$data[0,7]  = "fra"
$data[0,8]  = $true
$data[0,9]  = "superadmin"
$data[0,10] = $crDtimes
$data[0,11] = "NULL"
$data[0,12] = "NULL"
$data[0,13] = $false
$data[0,14] = "NULL"

# Row 3 (id 327 - High Speed Dual Iris Scanner / Eyecool)
$data[1,0]  = 327
$data[1,1]  = "High Speed Dual Iris Scanner"
$data[1,2]  = "Eyecool"
$data[1,3]  = "ECI102"
$data[1,4]  = "IRS"
$data[1,5]  = 2.34
$data[1,6]  = "Pour capturer l'iris"
$data[1,7]  = "fra"
$data[1,8]  = $true
$data[1,9]  = "superadmin"
$data[1,10] = $crDtimes
$data[1,11] = "NULL"
$data[1,12] = "NULL"
$data[1,13] = $false
$data[1,14] = "NULL"

# Row 4 (id 736 - Webcam / Logitech)
$data[2,0]  = 736
$data[2,1]  = "Webcam"
$data[2,2]  = "Logitech"
$data[2,3]  = "C270"
$data[2,4]  = "CMR"
$data[2,5]  = 2.086
$data[2,6]  = "Pour prendre la photo"
$data[2,7]  = "fra"
$data[2,8]  = $true
$data[2,9]  = "superadmin"
$data[2,10] = $crDtimes
$data[2,11] = "NULL"
$data[2,12] = "NULL"
$data[2,13] = $false
$data[2,14] = "NULL"

# Row 5 (id 801 - imageFORMULA / Canon) -- descr kept verbatim with the
# mangled accent characters exactly as found in the source data.
$data[3,0]  = 801
$data[3,1]  = "imageFORMULA"
$data[3,2]  = "Canon"
$data[3,3]  = "DR-C130"
$data[3,4]  = "SCN"
$data[3,5]  = 1.02
$data[3,6]  = "Pour numÃ©riser les documents"
$data[3,7]  = "fra"
$data[3,8]  = $true
$data[3,9]  = "superadmin"
$data[3,10] = $crDtimes
$data[3,11] = "NULL"
$data[3,12] = "NULL"
$data[3,13] = $false
$data[3,14] = "NULL"

# Row 6 (id 920 - Single Function Inkjet / Canon)
$data[4,0]  = 920
$data[4,1]  = "Single Function Inkjet"
$data[4,2]  = "Canon"
$data[4,3]  = "TS207"
$data[4,4]  = "PRT"
$data[4,5]  = 1.123
$data[4,6]  = "Pour imprimer les documents"
$data[4,7]  = "fra"
$data[4,8]  = $true
$data[4,9]  = "superadmin"
$data[4,10] = $crDtimes
$data[4,11] = "NULL"
$data[4,12] = "NULL"
$data[4,13] = $false
$data[4,14] = "NULL"

# Row 7 (id 444 - Fingerprint Scanner / Safran Morpho)
$data[5,0]  = 444
$data[5,1]  = "Fingerprint Scanner"
$data[5,2]  = "Safran Morpho"
$data[5,3]  = "1300 E2"
$data[5,4]  = "FRS"
$data[5,5]  = 1.12
$data[5,6]  = "Pour capturer les empreintes digitales"
$data[5,7]  = "fra"
$data[5,8]  = $true
$data[5,9]  = "superadmin"
$data[5,10] = $crDtimes
$data[5,11] = "NULL"
$data[5,12] = "NULL"
$data[5,13] = $false
$data[5,14] = "NULL"

$ws.Range("A2:O7").Value = $data

# The cr_dtimes column (K) picks up a date/time number format, same as
# Excel auto-applies when a date-like serial number is entered.
$ws.Range("K2:K7").NumberFormat = "mm:ss.0"

# Restore the active selection reported for the edited sheet.
$ws.Range("E17").Select()
